$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1234603333333333
$ws.Range("H2").Value = 0.370381
$ws.Range("I2").Value = 0.002558470358543426
$ws.Range("J2").Value = 0.002636284444771545
$ws.Range("M2").Value = 2.005664333333333
$ws.Range("N2").Value = 6.016992999999999
$ws.Range("O2").Value = 0.02976781902817159
$ws.Range("P2").Value = 0.03172257287647481
$ws.Range("Q2").Value = 0.2476199871481111
$ws.Range("R2").Value = 2.228579884333
$ws.Range("S2").Value = 0.00007616008262206199
$ws.Range("T2").Value = 0.00008362972542238229

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1234603333333333
$ws.Range("H3").Value = 0.370381
$ws.Range("I3").Value = 0.002558470358543426
$ws.Range("J3").Value = 0.002636284444771545
$ws.Range("O3").Value = 0.4846964599741412
$ws.Range("P3").Value = 0.5165248673390457
$ws.Range("Q3").Value = 4.031888633693556
$ws.Range("R3").Value = 36.286997703242
$ws.Range("S3").Value = 0.00124008152573477
$ws.Range("T3").Value = 0.001361706473103612

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1234603333333333
$ws.Range("H4").Value = 0.370381
$ws.Range("I4").Value = 0.002558470358543426
$ws.Range("J4").Value = 0.002636284444771545
$ws.Range("M4").Value = 10.495667
$ws.Range("N4").Value = 31.487001
$ws.Range("O4").Value = 0.1557753760903259
$ws.Range("P4").Value = 0.1660046278737794
$ws.Range("Q4").Value = 1.295798546375667
$ws.Range("R4").Value = 11.662186917381
$ws.Range("S4").Value = 0.000398546682318053
$ws.Range("T4").Value = 0.0004376354182237335

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1234603333333333
$ws.Range("H5").Value = 0.370381
$ws.Range("I5").Value = 0.002558470358543426
$ws.Range("J5").Value = 0.002636284444771545
$ws.Range("M5").Value = 12.4553565
$ws.Range("N5").Value = 24.910713
$ws.Range("O5").Value = 0.1848608423958749
$ws.Range("P5").Value = 0.1313333601264699
$ws.Range("Q5").Value = 1.5377424652755
$ws.Range("R5").Value = 9.226454791653001
$ws.Range("S5").Value = 0.0004729609857252137
$ws.Range("T5").Value = 0.0003462320943809921

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1234603333333333
$ws.Range("H6").Value = 0.370381
$ws.Range("I6").Value = 0.002558470358543426
$ws.Range("J6").Value = 0.002636284444771545
$ws.Range("M6").Value = 9.762884
$ws.Range("N6").Value = 29.288652
$ws.Range("O6").Value = 0.1448995025114864
$ws.Range("P6").Value = 0.1544145717842301
$ws.Range("Q6").Value = 1.205328912934667
$ws.Range("R6").Value = 10.847960216412
$ws.Range("S6").Value = 0.0003707210821433266
$ws.Range("T6").Value = 0.0004070807336408249

# Row 7
$ws.Range("I7").Value = 0.9088918061291337
$ws.Range("J7").Value = 0.9365351146153917
$ws.Range("M7").Value = 2.005664333333333
$ws.Range("N7").Value = 6.016992999999999
$ws.Range("O7").Value = 0.02976781902817159
$ws.Range("P7").Value = 0.03172257287647481
$ws.Range("Q7").Value = 87.96653695876677
$ws.Range("R7").Value = 791.698832628901
$ws.Range("S7").Value = 0.02705572680104007
$ws.Range("T7").Value = 0.02970930342476446

# Row 8
$ws.Range("I8").Value = 0.9088918061291337
$ws.Range("J8").Value = 0.9365351146153917
$ws.Range("O8").Value = 0.4846964599741412
$ws.Range("P8").Value = 0.5165248673390457
$ws.Range("S8").Value = 0.4405366409302946
$ws.Range("T8").Value = 0.4837436758350732

# Row 9
$ws.Range("I9").Value = 0.9088918061291337
$ws.Range("J9").Value = 0.9365351146153917
$ws.Range("M9").Value = 10.495667
$ws.Range("N9").Value = 31.487001
$ws.Range("O9").Value = 0.1557753760903259
$ws.Range("P9").Value = 0.1660046278737794
$ws.Range("Q9").Value = 460.3300082262397
$ws.Range("R9").Value = 4142.970074036158
$ws.Range("S9").Value = 0.1415829629251813
$ws.Range("T9").Value = 0.1554691631924554

# Row 10
$ws.Range("I10").Value = 0.9088918061291337
$ws.Range("J10").Value = 0.9365351146153917
$ws.Range("M10").Value = 12.4553565
$ws.Range("N10").Value = 24.910713
$ws.Range("O10").Value = 0.1848608423958749
$ws.Range("P10").Value = 0.1313333601264699
$ws.Range("Q10").Value = 546.2801325638235
$ws.Range("R10").Value = 3277.680795382941
$ws.Range("S10").Value = 0.1680185049277398
$ws.Range("T10").Value = 0.122998303478868

# Row 11
$ws.Range("I11").Value = 0.9088918061291337
$ws.Range("J11").Value = 0.9365351146153917
$ws.Range("M11").Value = 9.762884
$ws.Range("N11").Value = 29.288652
$ws.Range("O11").Value = 0.1448995025114864
$ws.Range("P11").Value = 0.1544145717842301
$ws.Range("Q11").Value = 428.1908402802626
$ws.Range("R11").Value = 3853.717562522364
$ws.Range("S11").Value = 0.1316979705448778
$ws.Range("T11").Value = 0.1446146686842305

# Row 12
$ws.Range("G12").Value = 4.2730135
$ws.Range("H12").Value = 8.546027
$ws.Range("I12").Value = 0.08854972351232299
$ws.Range("J12").Value = 0.06082860093983664
$ws.Range("M12").Value = 2.005664333333333
$ws.Range("N12").Value = 6.016992999999999
$ws.Range("O12").Value = 0.02976781902817159
$ws.Range("P12").Value = 0.03172257287647481
$ws.Range("Q12").Value = 8.570230772801834
$ws.Range("R12").Value = 51.421384636811
$ws.Range("S12").Value = 0.002635932144509461
$ws.Range("T12").Value = 0.001929639726287972

# Row 13
$ws.Range("G13").Value = 4.2730135
$ws.Range("H13").Value = 8.546027
$ws.Range("I13").Value = 0.08854972351232299
$ws.Range("J13").Value = 0.06082860093983664
$ws.Range("O13").Value = 0.4846964599741412
$ws.Range("P13").Value = 0.5165248673390457
$ws.Range("Q13").Value = 139.5453430030357
$ws.Range("R13").Value = 837.272058018214
$ws.Range("S13").Value = 0.04291973751811193
$ws.Range("T13").Value = 0.03141948503086887

# Row 14
$ws.Range("G14").Value = 4.2730135
$ws.Range("H14").Value = 8.546027
$ws.Range("I14").Value = 0.08854972351232299
$ws.Range("J14").Value = 0.06082860093983664
$ws.Range("M14").Value = 10.495667
$ws.Range("N14").Value = 31.487001
$ws.Range("O14").Value = 0.1557753760903259
$ws.Range("P14").Value = 0.1660046278737794
$ws.Range("Q14").Value = 44.84812678250451
$ws.Range("R14").Value = 269.088760695027
$ws.Range("S14").Value = 0.01379386648282648
$ws.Range("T14").Value = 0.01009782926310021

# Row 15
$ws.Range("G15").Value = 4.2730135
$ws.Range("H15").Value = 8.546027
$ws.Range("I15").Value = 0.08854972351232299
$ws.Range("J15").Value = 0.06082860093983664
$ws.Range("M15").Value = 12.4553565
$ws.Range("N15").Value = 24.910713
$ws.Range("O15").Value = 0.1848608423958749
$ws.Range("P15").Value = 0.1313333601264699
$ws.Range("Q15").Value = 53.22190647181276
$ws.Range("R15").Value = 212.887625887251
$ws.Range("S15").Value = 0.01636937648240983
$ws.Range("T15").Value = 0.00798882455322089

# Row 16
$ws.Range("G16").Value = 4.2730135
$ws.Range("H16").Value = 8.546027
$ws.Range("I16").Value = 0.08854972351232299
$ws.Range("J16").Value = 0.06082860093983664
$ws.Range("M16").Value = 9.762884
$ws.Range("N16").Value = 29.288652
$ws.Range("O16").Value = 0.1448995025114864
$ws.Range("P16").Value = 0.1544145717842301
$ws.Range("Q16").Value = 41.716935130934
$ws.Range("R16").Value = 250.301610785604
$ws.Range("S16").Value = 0.01283081088446527
$ws.Range("T16").Value = 0.009392822366358689
